$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet as the last tab ---
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Copy the header / date-column formatting from "Weekly Quantity" so the new
# sheet matches the workbook's existing look: bold + centered + bordered
# header row, and a date number format running down column A.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2:A12").Copy()
$wsForecast.Range("A2:A20").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast rows: ds (order-week date serial), PO_Forecast, yhat_lower, yhat_upper
$data = @(
  @(44948.99999999999, 32, 3.948480590707766, 61.57774640290018),
  @(44955.99999999999, 31, 2.317529242613807, 61.46803324112675),
  @(44962.99999999999, 29, -1.653827918756607, 59.04364420868994),
  @(44969.99999999999, 27, -3.944374171408937, 54.61105273295123),
  @(44976.99999999999, 25, -5.157595250351463, 56.25398409448614),
  @(44983.99999999999, 23, -5.299966855975468, 55.07869297573436),
  @(44990.99999999999, 21, -9.381587569624275, 50.40527625869287),
  @(44997.99999999999, 20, -11.80177947533241, 50.0276723406017),
  @(45004.99999999999, 18, -11.3919967450157, 48.14434777679232),
  @(45011.99999999999, 16, -13.97853927819847, 48.43111761811893),
  @(45018.99999999999, 14, -15.28303611252743, 44.60379218394732),
  @(45025.99999999999, 12, -17.94896658131541, 42.94091217601816),
  @(45032.99999999999, 11, -18.08540626933987, 40.94400984880694),
  @(45039.99999999999, 9, -17.3889537073075, 37.65242603042372),
  @(45046.99999999999, 7, -23.56824124016689, 36.50493102812035),
  @(45053.99999999999, 5, -26.7733374569099, 33.43865679952293),
  @(45060.99999999999, 3, -26.81831526273327, 34.617214098917),
  @(45067.99999999999, 1, -27.7831066946515, 30.41383177156771),
  @(45074.99999999999, 0, -30.64807556018309, 29.01467173966362)
)

$r = 2
foreach ($row in $data) {
  $wsForecast.Cells.Item($r, 1).Value = $row[0]
  $wsForecast.Cells.Item($r, 2).Value = $row[1]
  $wsForecast.Cells.Item($r, 3).Value = $row[2]
  $wsForecast.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}
